# Segunda parte de la clase del 25/05/2015:
# limpia los datos de ejemplo que se habian cargado para la demo de
# validacion de datos, deja las validaciones listas para que los alumnos
# las vuelvan a crear, y quita el defined name heredado del addin MySQL for
# Excel que ya no se usa.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Defined names: el workbook traia un nombre oculto del plugin MySQL for
#     Excel que ya no hace falta; dejamos solo CURSOS. ---
foreach ($n in @($wb.Names)) {
    if ($n.Name -eq "LOCAL_MYSQL_DATE_FORMAT") {
        $n.Delete()
    }
}

# --- Fila 3 (ALU-0001): se quitan los datos de ejemplo de sexo, telefono,
#     curso y las notas de practica que se habian puesto como demo. ---
$ws.Range("D3:F3").ClearContents()
$ws.Range("H3:K3").ClearContents()

# --- Fila 4 (ALU-0002): se quitan sexo y telefono de ejemplo. ---
$ws.Range("D4:E4").ClearContents()

# --- Fila 5: el promedio (columna L) ahora usa la misma formula compartida
#     que el resto de la columna (antes estaba vacia). ---
$ws.Range("L5").Formula = '=IFERROR(ROUND(AVERAGE(H5:K5),0),"")'

# --- Fila 6: se quita la fecha de nacimiento de ejemplo. ---
$ws.Range("G6").ClearContents()

# --- N10/N11 contenian la lista auxiliar Masculino/Femenino que ya no se
#     usa (las cadenas quedan huerfanas y se eliminan del workbook). ---
$ws.Range("N10").ClearContents()
$ws.Range("N11").ClearContents()

# --- Cursos de ejemplo que se habian asignado a algunos alumnos se quitan. ---
$ws.Range("F12").ClearContents()
$ws.Range("F17").ClearContents()
$ws.Range("F19").ClearContents()

# --- Columna F (CURSO) un poco mas angosta. ColumnWidth esta en "caracteres"
#     y Excel le suma un margen fijo (5px ~ 0.8333 caracteres con la fuente
#     por defecto) antes de guardar el ancho crudo en el XML, asi que restamos
#     ese margen para que el ancho guardado quede en exactamente 26. ---
$ws.Columns.Item(6).ColumnWidth = 26 - 0.8333333333333334

# --- Ya no se necesitan las validaciones de datos (se volveran a armar en
#     la siguiente sesion de clase). ---
$ws.Range("D3:D22").Validation.Delete()
$ws.Range("E3:E22").Validation.Delete()
$ws.Range("H3:L22").Validation.Delete()
$ws.Range("F3:F22").Validation.Delete()
$ws.Range("G3:G22").Validation.Delete()

# --- Selecciona L13 (sin scroll fijo en A2) para que quede igual a como el
#     profesor dejo el archivo al guardar. ---
$ws.Range("L13").Select() | Out-Null
